# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 26 (pushing the existing Pomelo
# price-history rows 26..51 down to 27..52) and populate the new row
# with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 26; Excel shifts rows 26:51
# down to 27:52 and extends the used range to A1:T52.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with this week's record.
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44469
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100102
$ws.Range("H26").Value = "Cítricos"
$ws.Range("I26").Value = 100102006
$ws.Range("J26").Value = "Pomelo"
$ws.Range("K26").Value = "Start Ruby"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 450
$ws.Range("N26").Value = 7000
$ws.Range("O26").Value = 7500
$ws.Range("P26").Value = 7278
$ws.Range("Q26").Value = "$/caja 14 kilos granel"
$ws.Range("R26").Value = "Región Metropolitana"
$ws.Range("S26").Value = 520
$ws.Range("T26").Value = 14
